# Saptamana 27 -> update: extend Sheet1 totals/summary column (H) and add a new
# Sheet2 with "Caracteristicile cazurilor confirmate si a deceselor".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$xlCenter = -4108

# ---------------------------------------------------------------------------
# Sheet1 changes
# ---------------------------------------------------------------------------

# New column widths: G narrower, H (new) wide for the comorbidity notes.
$ws1.Columns.Item(7).ColumnWidth = 12.7109375
$ws1.Columns.Item(8).ColumnWidth = 47.7109375

# Three note cells next to the first data rows, centered like the rest of
# the sheet.
$ws1.Range("H3").Value = "94% din decese aveau comorbiditati asociate"
$ws1.Range("H4").Value = "78.4% din decese au fost la persoane de peste 60 ani"
$ws1.Range("H5").Value = "59.2% din decese au fost la barbati"
$ws1.Range("H3:H5").HorizontalAlignment = $xlCenter

# Header row (B1:F1) keeps its centered style - re-assert alignment so the
# style table matches (Excel creates a fresh xf entry for this).
$ws1.Range("B1:F1").HorizontalAlignment = $xlCenter

# Extend the "Total" row (46) with the grand totals, and add a new "-" row
# (47) under it.
$ws1.Range("C46").Value = 32948
$ws1.Range("D46").Value = 1947
$ws1.Range("E46").Value = 3725
$ws1.Range("F46").Value = 134
$ws1.Range("C46:F46").HorizontalAlignment = $xlCenter

$ws1.Range("C47").Value = "-"
$ws1.Range("D47").Value = "-"
$ws1.Range("E47").Value = "-"
$ws1.Range("F47").Value = "-"

# Selection / view restore for Sheet1 (it will no longer be the active tab).
$ws1.Range("H3:H5").Select()

# ---------------------------------------------------------------------------
# New Sheet2: "Caracteristicile cazurilor confirmate si a deceselor"
# ---------------------------------------------------------------------------

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1:O1").Merge()
$ws2.Range("A1").Value = "Caracteristicile cazurilor confirmate si a deceselor"

$ws2.Range("B3:F3").Merge()
$ws2.Range("B3").Value = "Cazuri"
$ws2.Range("G3:P3").Merge()
$ws2.Range("G3").Value = "Decese"

$ws2.Range("B4").Value = "Varsta,mediana(range)"
$ws2.Range("C4").Value = "Sex,masculin"
$ws2.Range("D4").Value = "Import"
$ws2.Range("E4").Value = "Vindecati"
$ws2.Range("F4").Value = "Personal sanitar"
$ws2.Range("G4").Value = "Varsta,mediana(range)"
$ws2.Range("H4").Value = "Sex,masculin"
$ws2.Range("I4").Value = "Afectiuni cardiovasculare"
$ws2.Range("J4").Value = "Diabet"
$ws2.Range("K4").Value = "Afectiuni neurologice"
$ws2.Range("L4").Value = "Afectiuni renale"
$ws2.Range("M4").Value = "Obezitate"
$ws2.Range("N4").Value = "Afectiuni Pulmonare"
$ws2.Range("O4").Value = "Neoplasm"
$ws2.Range("P4").Value = "Altele"

$ws2.Range("A5").Value = "Numar"
$ws2.Range("B5").Value = "49(0-99)"
$ws2.Range("C5").Value = 14978
$ws2.Range("D5").Value = 803
$ws2.Range("E5").Value = 23552
$ws2.Range("F5").Value = 3422
$ws2.Range("G5").Value = "69(20-99)"
$ws2.Range("H5").Value = 1125
$ws2.Range("I5").Value = 1256
$ws2.Range("J5").Value = 594
$ws2.Range("K5").Value = 428
$ws2.Range("L5").Value = 377
$ws2.Range("M5").Value = 330
$ws2.Range("N5").Value = 340
$ws2.Range("O5").Value = 230
$ws2.Range("P5").Value = 370

$ws2.Range("A6").Value = "Procentaj"
$ws2.Range("C6").Value = 45.5
$ws2.Range("D6").Value = 2.4
$ws2.Range("E6").Value = 71.5
$ws2.Range("F6").Value = 10.4
$ws2.Range("H6").Value = 59.2
$ws2.Range("I6").Value = 66.1
$ws2.Range("J6").Value = 31.2
$ws2.Range("K6").Value = 22.5
$ws2.Range("L6").Value = 19.8
$ws2.Range("M6").Value = 17.4
$ws2.Range("N6").Value = 17.9
$ws2.Range("O6").Value = 12.1
$ws2.Range("P6").Value = 19.5

# Whole used range is centered, like the source sheet.
$ws2.Range("A1:P6").HorizontalAlignment = $xlCenter

# Column widths tuned as in the source file.
$ws2.Columns.Item(1).ColumnWidth = 15.5703125
$ws2.Columns.Item(2).ColumnWidth = 22.5703125
$ws2.Columns.Item(3).ColumnWidth = 15
$ws2.Columns.Item(5).ColumnWidth = 10.7109375
$ws2.Columns.Item(6).ColumnWidth = 16
$ws2.Columns.Item(7).ColumnWidth = 23.42578125
$ws2.Columns.Item(8).ColumnWidth = 14.28515625
$ws2.Columns.Item(9).ColumnWidth = 24.5703125
$ws2.Columns.Item(11).ColumnWidth = 21.28515625
$ws2.Columns.Item(12).ColumnWidth = 17.140625
$ws2.Columns.Item(14).ColumnWidth = 20.42578125
$ws2.Columns.Item(15).ColumnWidth = 10.42578125

$ws2.Range("G7").Select()
$ws2.Activate()
